# Natmi following Dr Hou advice
#
# Updates the LR-pairs_lrc2p "Bsg-Slc16a7" result table:
#  - recomputed values for the existing 12 data rows (rows 2-13)
#  - appended 4 new data rows (rows 14-17) for the previously-missing
#    "sCs" sending-cluster block, extending the table from A1:T13 to A1:T17
#
# Column layout (A..T):
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E Ligand-expressing cells, F Ligand detection rate,
#   G Ligand average expression value, H Ligand total expression value,
#   I Ligand derived specificity of average expression value,
#   J Ligand derived specificity of total expression value,
#   K Receptor-expressing cells, L Receptor detection rate,
#   M Receptor average expression value, N Receptor total expression value,
#   O Receptor derived specificity of average expression value,
#   P Receptor derived specificity of total expression value,
#   Q Edge average expression weight, R Edge total expression weight,
#   S Edge average expression derived specificity,
#   T Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Bsg", "Slc16a7", "ECs", 3, 1, 36.78646733333333, 110.359402, 0.1576941973553631, 0.1576941973553631, 1, 0.3333333333333333, 1.158627, 3.475881, 0.09636025486884757, 0.09636025486884758, 42.62179428701801, 383.5961485831621, 0.01519545304850114, 0.01519545304850114),
  @("ECs", "Bsg", "Slc16a7", "FAPs", 3, 1, 36.78646733333333, 110.359402, 0.1576941973553631, 0.1576941973553631, 3, 1, 0.450595, 1.351785, 0.03747491560496033, 0.03747491560496033, 16.57579824806334, 149.18218423257, 0.005909576737284191, 0.005909576737284192),
  @("ECs", "Bsg", "Slc16a7", "M2", 3, 1, 36.78646733333333, 110.359402, 0.1576941973553631, 0.1576941973553631, 3, 1, 9.934199999999999, 29.8026, 0.8262038118549847, 0.8262038118549848, 365.4441237828, 3288.9971140452, 0.1302875469624132, 0.1302875469624133),
  @("ECs", "Bsg", "Slc16a7", "sCs", 3, 1, 36.78646733333333, 110.359402, 0.1576941973553631, 0.1576941973553631, 3, 1, 0.4804876666666666, 1.441463, 0.03996101767120728, 0.03996101767120728, 17.67544385390289, 159.078994685126, 0.006301620607164513, 0.006301620607164515),
  @("FAPs", "Bsg", "Slc16a7", "ECs", 3, 1, 103.9192913333333, 311.757874, 0.445475481188675, 0.4454754811886751, 1, 0.3333333333333333, 1.158627, 3.475881, 0.09636025486884757, 0.09636025486884758, 120.403696759666, 1083.633270836994, 0.04292613090516324, 0.04292613090516325),
  @("FAPs", "Bsg", "Slc16a7", "FAPs", 3, 1, 103.9192913333333, 311.757874, 0.445475481188675, 0.4454754811886751, 3, 1, 0.450595, 1.351785, 0.03747491560496033, 0.03747491560496033, 46.82551307834333, 421.42961770509, 0.01669415606162469, 0.01669415606162469),
  @("FAPs", "Bsg", "Slc16a7", "M2", 3, 1, 103.9192913333333, 311.757874, 0.445475481188675, 0.4454754811886751, 3, 1, 9.934199999999999, 29.8026, 0.8262038118549847, 0.8262038118549848, 1032.3550239636, 9291.1952156724, 0.3680535406460169, 0.3680535406460169),
  @("FAPs", "Bsg", "Slc16a7", "sCs", 3, 1, 103.9192913333333, 311.757874, 0.445475481188675, 0.4454754811886751, 3, 1, 0.4804876666666666, 1.441463, 0.03996101767120728, 0.03996101767120728, 49.93193781440689, 449.387440329662, 0.01780165357587021, 0.01780165357587021),
  @("M2", "Bsg", "Slc16a7", "ECs", 3, 1, 43.05432033333333, 129.162961, 0.1845628835768525, 0.1845628835768525, 1, 0.3333333333333333, 1.158627, 3.475881, 0.09636025486884757, 0.09636025486884758, 49.883898004849, 448.955082043641, 0.01778452650079495, 0.01778452650079496),
  @("M2", "Bsg", "Slc16a7", "FAPs", 3, 1, 43.05432033333333, 129.162961, 0.1845628835768525, 0.1845628835768525, 3, 1, 0.450595, 1.351785, 0.03747491560496033, 0.03747491560496033, 19.40006147059833, 174.600553235385, 0.006916478485850667, 0.006916478485850668),
  @("M2", "Bsg", "Slc16a7", "M2", 3, 1, 43.05432033333333, 129.162961, 0.1845628835768525, 0.1845628835768525, 3, 1, 9.934199999999999, 29.8026, 0.8262038118549847, 0.8262038118549848, 427.7102290553999, 3849.3920614986, 0.1524865579381433, 0.1524865579381434),
  @("M2", "Bsg", "Slc16a7", "sCs", 3, 1, 43.05432033333333, 129.162961, 0.1845628835768525, 0.1845628835768525, 3, 1, 0.4804876666666666, 1.441463, 0.03996101767120728, 0.03996101767120728, 20.68706991688255, 186.183629251943, 0.007375320652063575, 0.007375320652063576),
  @("sCs", "Bsg", "Slc16a7", "ECs", 3, 1, 49.51716233333334, 148.551487, 0.2122674378791094, 0.2122674378791094, 1, 0.3333333333333333, 1.158627, 3.475881, 0.09636025486884757, 0.09636025486884758, 57.37192124278301, 516.3472911850471, 0.02045414441438825, 0.02045414441438825),
  @("sCs", "Bsg", "Slc16a7", "FAPs", 3, 1, 49.51716233333334, 148.551487, 0.2122674378791094, 0.2122674378791094, 3, 1, 0.450595, 1.351785, 0.03747491560496033, 0.03747491560496033, 22.31218576158834, 200.809671854295, 0.007954704320200782, 0.007954704320200784),
  @("sCs", "Bsg", "Slc16a7", "M2", 3, 1, 49.51716233333334, 148.551487, 0.2122674378791094, 0.2122674378791094, 3, 1, 9.934199999999999, 29.8026, 0.8262038118549847, 0.8262038118549848, 491.9133940518, 4427.2205464662, 0.1753761663084113, 0.1753761663084114),
  @("sCs", "Bsg", "Slc16a7", "sCs", 3, 1, 49.51716233333334, 148.551487, 0.2122674378791094, 0.2122674378791094, 3, 1, 0.4804876666666666, 1.441463, 0.03996101767120728, 0.03996101767120728, 23.79238578949789, 214.131472105481, 0.008482422836108982, 0.008482422836108984)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $rowArr = $data[$i]
  $r = $startRow + $i
  for ($c = 0; $c -lt $rowArr.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowArr[$c]
  }
}

Write-Host "Updated $($data.Count) data rows (A2:T$($startRow + $data.Count - 1))"
